$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Regular_Timetable")
$ws.Range("B3").Value = 'MA262 [C003]'
$ws.Range("C3").Value = 'MA262 [C003]'
$ws.Range("D3").Value = 'DA262 [C304]'
$ws.Range("E3").Value = 'CS304 [C205]'
$ws.Range("F3").Value = 'CS307 [C302]'
$ws.Range("B4").Value = 'MA261 [C003]'
$ws.Range("C4").Value = 'MA261 [C003]'
$ws.Range("D4").Value = 'CS304 [C205]'
$ws.Range("E4").Value = 'CS307 [C302]'
$ws.Range("C6").Value = 'DA262 [C304]'
$ws.Range("E6").Value = 'CS307 (Lab) [L207]'
$ws.Range("B7").Value = 'CS304 (Tutorial) [C304]'
$ws.Range("E7").Value = 'CS307 (Lab) [L207]'
$ws.Range("B8").Value = 'DA261 [C204]'
$ws.Range("C8").Value = 'DA261 (Lab) [L106]'
$ws.Range("D8").Value = 'DA262 (Lab) [L207]'
$ws.Range("C9").Value = 'DA261 (Lab) [L106]'
$ws.Range("D9").Value = 'DA262 (Lab) [L207]'
$ws.Range("D25").Value = 'Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]'
$ws.Range("E25").Value = 'Tue 14:30-15:30 [C101]'
$ws.Range("D26").Value = 'Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]'
$ws.Range("E26").Value = 'Tue 14:30-15:30 [C102]'
$ws.Range("E27").Value = 'Tue 14:30-15:30 [C104]'
$ws.Range("D28").Value = 'Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]'
$ws.Range("E28").Value = 'Tue 14:30-15:30 [C202]'
$ws.Range("D29").Value = 'Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]'
$ws.Range("E29").Value = 'Tue 14:30-15:30 [C203]'
$ws.Range("D30").Value = 'Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]'
$ws.Range("E30").Value = 'Tue 14:30-15:30 [C204]'
$ws.Range("D31").Value = 'Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]'
$ws.Range("E31").Value = 'Tue 14:30-15:30 [C205]'

$ws = $wb.Worksheets.Item("PreMid_Timetable")
$ws.Range("D3").Value = 'DA262 [C305]'
$ws.Range("E3").Value = 'CS304 [C101]'
$ws.Range("D4").Value = 'CS304 [C302]'
$ws.Range("E4").Value = 'CS307 [C303]'
$ws.Range("C6").Value = 'DA262 [C305]'
$ws.Range("E6").Value = 'CS307 (Lab) [L207]'
$ws.Range("B7").Value = 'CS304 (Tutorial) [C305]'
$ws.Range("E7").Value = 'CS307 (Lab) [L207]'
$ws.Range("C8").Value = 'DA262 (Lab) [L106]'
$ws.Range("D8").Value = 'CS307 [C302]'
$ws.Range("C9").Value = 'DA262 (Lab) [L106]'
$ws.Range("D25").Value = 'Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]'
$ws.Range("E25").Value = 'Tue 14:30-15:30 [C101]'
$ws.Range("D26").Value = 'Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]'
$ws.Range("E26").Value = 'Tue 14:30-15:30 [C102]'
$ws.Range("E27").Value = 'Tue 14:30-15:30 [C104]'
$ws.Range("D28").Value = 'Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]'
$ws.Range("E28").Value = 'Tue 14:30-15:30 [C202]'
$ws.Range("D29").Value = 'Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]'
$ws.Range("E29").Value = 'Tue 14:30-15:30 [C203]'
$ws.Range("D30").Value = 'Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]'
$ws.Range("E30").Value = 'Tue 14:30-15:30 [C204]'
$ws.Range("D31").Value = 'Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]'
$ws.Range("E31").Value = 'Tue 14:30-15:30 [C205]'

$ws = $wb.Worksheets.Item("PostMid_Timetable")
$ws.Range("D3").Value = 'DA262 [C102]'
$ws.Range("E3").Value = 'CS304 [C102]'
$ws.Range("D4").Value = 'CS304 [C303]'
$ws.Range("E4").Value = 'CS307 [C304]'
$ws.Range("C6").Value = 'DA262 [C102]'
$ws.Range("E6").Value = 'CS307 (Lab) [L207]'
$ws.Range("B7").Value = 'CS304 (Tutorial) [C102]'
$ws.Range("E7").Value = 'CS307 (Lab) [L207]'
$ws.Range("C8").Value = 'DA262 (Lab) [L106]'
$ws.Range("D8").Value = 'CS307 [C303]'
$ws.Range("E8").Value = 'DA261 [C204]'
$ws.Range("F8").Value = 'DA261 (Lab) [L106]'
$ws.Range("C9").Value = 'DA262 (Lab) [L106]'
$ws.Range("F9").Value = 'DA261 (Lab) [L106]'
$ws.Range("D25").Value = 'Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]'
$ws.Range("E25").Value = 'Tue 14:30-15:30 [C101]'
$ws.Range("D26").Value = 'Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]'
$ws.Range("E26").Value = 'Tue 14:30-15:30 [C102]'
$ws.Range("E27").Value = 'Tue 14:30-15:30 [C104]'
$ws.Range("D28").Value = 'Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]'
$ws.Range("E28").Value = 'Tue 14:30-15:30 [C202]'
$ws.Range("D29").Value = 'Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]'
$ws.Range("E29").Value = 'Tue 14:30-15:30 [C203]'
$ws.Range("D30").Value = 'Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]'
$ws.Range("E30").Value = 'Tue 14:30-15:30 [C204]'
$ws.Range("D31").Value = 'Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]'
$ws.Range("E31").Value = 'Tue 14:30-15:30 [C205]'

